$d = $word.ActiveDocument

# --- Update the title date line ("2024-11-24 Sunday" -> "2024-11-25 Monday") ---
$d.Content.Find.Execute("2024-11-24 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-11-25 Monday", 2)

# --- Update every arithmetic-problem cell in the practice table ---
# The new answer key, in document order (row-major: row 1 col 1..5, row 2 col 1..5, ...).
# Built this way (rather than global Find/Replace) because several of the
# original expressions repeat verbatim at different table positions but map
# to different replacement text, and some replacement texts collide with
# each other too, so a positional walk is required for correctness.
$values = @(
    "21+59=80",
    "47+31=78",
    "45-19=26",
    "95-32=63",
    "9+1=10",
    "42-3=39",
    "40+51=91",
    "91-4=87",
    "45+3=48",
    "8+28=36",
    "57+13=70",
    "19-8=11",
    "2+55=57",
    "72-68=4",
    "97-7=90",
    "28+3=31",
    "0+12=12",
    "45-27=18",
    "89-30=59",
    "67+21=88",
    "35+44=79",
    "73+22=95",
    "66-62=4",
    "52-15=37",
    "16+4=20",
    "20+15=35",
    "33+11=44",
    "60+0=60",
    "3+15=18",
    "87-80=7",
    "65+14=79",
    "86-71=15",
    "40-38=2",
    "43+21=64",
    "22-5=17",
    "49+3=52",
    "15+58=73",
    "3+91=94",
    "59-0=59",
    "5+54=59",
    "96-65=31",
    "52-31=21",
    "19-1=18",
    "3+22=25",
    "40-0=40",
    "77-63=14",
    "45-22=23",
    "56-41=15",
    "78-61=17",
    "73-10=63",
    "94+1=95",
    "84-71=13",
    "81-64=17",
    "68+20=88",
    "56+8=64",
    "82+3=85",
    "56+19=75",
    "18+16=34",
    "97-72=25",
    "62+18=80",
    "68-65=3",
    "17+27=44",
    "73+18=91",
    "17+44=61",
    "23+3=26",
    "70+8=78",
    "79-52=27",
    "24-3=21",
    "19+9=28",
    "37+28=65",
    "5+33=38",
    "21+59=80",
    "58+14=72",
    "81-72=9",
    "56-44=12",
    "69-6=63",
    "96+0=96",
    "62-3=59",
    "82-4=78",
    "22+16=38",
    "91-76=15",
    "31+41=72",
    "56+38=94",
    "86-65=21",
    "28+30=58",
    "61+3=64",
    "8+18=26",
    "88-36=52",
    "80-15=65",
    "96-91=5",
    "83-16=67",
    "33-32=1",
    "73-55=18",
    "63-3=60",
    "67+6=73",
    "44+35=79",
    "6+79=85",
    "59-8=51",
    "64+30=94",
    "17+7=24"
)


$tbl = $d.Tables.Item(1)
$rows = $tbl.Rows
$idx = 0
for ($r = 1; $r -le $rows.Count; $r++) {
    $cells = $rows.Item($r).Cells
    for ($c = 1; $c -le $cells.Count; $c++) {
        $cell = $cells.Item($c)
        $cell.Range.Text = $values[$idx]
        $idx = $idx + 1
    }
}

Write-Host "Updated $idx table cells."
